$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/12/2024  Through  2/18/2024"

# Row 16
$src = $ws.Range("D16")
$dst = $ws.Range("C16")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$dst.Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 25
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -16.666666666666
$ws.Range("N16").Value = -79.591836734693

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -5.263157894736
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 29
$ws.Range("K17").Value = -20.689655172413
$ws.Range("L17").Value = 109.090909090909
$ws.Range("M17").Value = 4.545454545454
$ws.Range("N17").Value = -28.125

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 13
$ws.Range("K18").Value = -23.076923076923
$ws.Range("L18").Value = 100
$ws.Range("M18").Value = -60
$ws.Range("N18").Value = -93.865030674846

# Row 19
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 53
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 60.60606060606
$ws.Range("I19").Value = 83
$ws.Range("J19").Value = 71
$ws.Range("K19").Value = 16.901408450704
$ws.Range("L19").Value = 84.444444444444
$ws.Range("M19").Value = 50.90909090909
$ws.Range("N19").Value = -29.661016949152

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = -55.555555555555
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -58.823529411764
$ws.Range("L20").Value = -58.823529411764
$ws.Range("M20").Value = -22.222222222222
$ws.Range("N20").Value = -98.172323759791

# Row 21
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 87
$ws.Range("G21").Value = 73
$ws.Range("H21").Value = 19.17808219178
$ws.Range("I21").Value = 133
$ws.Range("J21").Value = 138
$ws.Range("K21").Value = -3.623188405797
$ws.Range("L21").Value = 58.333333333333
$ws.Range("M21").Value = 5.555555555555
$ws.Range("N21").Value = -82.195448460508

# Row 23
$ws.Range("D23").Value = 1
$ws.Range("G23").Value = 5
$ws.Range("J23").Value = 6

# Row 24
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 5.555555555555
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 77
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 132
$ws.Range("J24").Value = 132
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 40.425531914893
$ws.Range("M24").Value = -36.538461538461

# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 350
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 107.142857142857
$ws.Range("I25").Value = 39
$ws.Range("J25").Value = 28
$ws.Range("K25").Value = 39.285714285714
$ws.Range("L25").Value = 62.5
$ws.Range("M25").Value = -30.357142857142

# Row 27
$ws.Range("C27").Value = 1
$src = $ws.Range("C27")
$dst = $ws.Range("D27")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$dst.Value = 2
$src = $ws.Range("H27")
$dst = $ws.Range("E27")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$dst.Value = -50
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 4
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = -33.333333333333

# Row 30
$src = $ws.Range("C30")
$dst = $ws.Range("D30")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$src.Copy()
$dst.PasteSpecial(-4163)
$excel.CutCopyMode = $false
$src = $ws.Range("M26")
$dst = $ws.Range("E30")
$src.Copy()
$dst.PasteSpecial(-4122)
$excel.CutCopyMode = $false
$src.Copy()
$dst.PasteSpecial(-4163)
$excel.CutCopyMode = $false

